$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6704119850187266
$ws1.Range("C2").Value = 0.6117936117936118
$ws1.Range("D2").Value = 0.9325842696629213
$ws1.Range("E2").Value = 0.7388724035608308
$ws1.Range("F2").Value = 0.8440677966101695
$ws1.Range("G2").Value = 0.9141485456085852
$ws1.Range("H2").Value = 0.7866360869138297
$ws1.Range("I2").Value = 498
$ws1.Range("J2").Value = 316
$ws1.Range("K2").Value = 218
$ws1.Range("L2").Value = 36

# --- Sheet "Classification Report" ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.8582677165354331
$ws2.Range("C2").Value = 0.4082397003745318
$ws2.Range("D2").Value = 0.5532994923857868

$ws2.Range("B3").Value = 0.6117936117936118
$ws2.Range("C3").Value = 0.9325842696629213
$ws2.Range("D3").Value = 0.7388724035608308

$ws2.Range("B4").Value = 0.6704119850187266
$ws2.Range("C4").Value = 0.6704119850187266
$ws2.Range("D4").Value = 0.6704119850187266
$ws2.Range("E4").Value = 0.6704119850187266

$ws2.Range("B5").Value = 0.7350306641645225
$ws2.Range("C5").Value = 0.6704119850187266
$ws2.Range("D5").Value = 0.6460859479733088

$ws2.Range("B6").Value = 0.7350306641645225
$ws2.Range("C6").Value = 0.6704119850187266
$ws2.Range("D6").Value = 0.6460859479733089

# --- Sheet "Confusion Matrix" ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 218
$ws3.Range("C2").Value = 316
$ws3.Range("B3").Value = 36
$ws3.Range("C3").Value = 498
